# Update TPM-derived values on the active sheet (rows 2-4) to reflect
# the new TPM calculation used by the data-generation scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.01249966666666667
$ws.Range("H2").Value = 0.037499
$ws.Range("M2").Value = 91.51130433333333
$ws.Range("N2").Value = 274.533913
$ws.Range("O2").Value = 0.9685519820468944
$ws.Range("P2").Value = 0.9685519820468945
$ws.Range("Q2").Value = 1.143860800398555
$ws.Range("R2").Value = 10.294747203587
$ws.Range("S2").Value = 0.9685519820468944
$ws.Range("T2").Value = 0.9685519820468945

# Row 3
$ws.Range("G3").Value = 0.01249966666666667
$ws.Range("H3").Value = 0.037499
$ws.Range("O3").Value = 0.001425786415744213
$ws.Range("P3").Value = 0.001425786415744214
$ws.Range("Q3").Value = 0.001683855096
$ws.Range("R3").Value = 0.015154695864
$ws.Range("S3").Value = 0.001425786415744213
$ws.Range("T3").Value = 0.001425786415744214

# Row 4
$ws.Range("G4").Value = 0.01249966666666667
$ws.Range("H4").Value = 0.037499
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509734999999999
$ws.Range("O4").Value = 0.03002223153736139
$ws.Range("P4").Value = 0.03002223153736139
$ws.Range("Q4").Value = 0.03545628364055555
$ws.Range("R4").Value = 0.3191065527649999
$ws.Range("S4").Value = 0.03002223153736139
$ws.Range("T4").Value = 0.03002223153736139
